$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range('D2').Value = '43.033.64'
$ws.Range('E2').Value = '  +0.10%  '

# Row 3: Ethereum
$ws.Range('D3').Value = '2.305.06'
$ws.Range('E3').Value = '  +0.12%  '

# Row 4: TetherUSD
$ws.Range('E4').Value = '  -0.05%  '

# Row 5: BNB
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '299.98'
$ws.Range('E5').Value = '  -0.65%  '

# Row 6: Solana
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '97.73'
$ws.Range('E6').Value = '  -0.75%  '

# Row 7: XRP
$ws.Range('E7').Value = '  -1.97%  '

# Row 8: USDC
$ws.Range('E8').Value = '  -0.06%  '

# Row 9: Cardano
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.507'
$ws.Range('E9').Value = '  -2.83%  '

# Row 10: Avalanche
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '35.86'
$ws.Range('E10').Value = '  +0.68%  '

# Row 11: Dogecoin
$ws.Range('E11').Value = '  +0.09%  '

# Row 12: Chainlink
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '18.16'
$ws.Range('E12').Value = '  +1.26%  '

# Row 13: TRON
$ws.Range('E13').Value = '  +1.78%  '

# Row 14: Polkadot
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.81'
$ws.Range('E14').Value = '  -1.03%  '

# Row 15: WrappedliquidstakedEther2.0
$ws.Range('D15').Value = '2.664.24'
$ws.Range('E15').Value = '  +0.04%  '

# Row 16: WrappedEther
$ws.Range('D16').Value = '2.296.87'
$ws.Range('E16').Value = '  -0.72%  '

# Row 17: Polygon
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.781'
$ws.Range('E17').Value = '  -0.98%  '

# Row 18: WrappedBTC
$ws.Range('D18').Value = '42.964.76'
$ws.Range('E18').Value = '  +0.15%  '

# Row 19: InternetComputer(DFINITY)
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.73'
$ws.Range('E19').Value = '  -4.95%  '

# Row 20: ShibaInu
$ws.Range('D20').Value = '0.0₃0904'
$ws.Range('E20').Value = '  -0.46%  '

# Row 21: Uniswap
$ws.Range('E21').Value = '  -1.91%  '

# Row 23: BitcoinCash
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '240.83'
$ws.Range('E23').Value = '  +0.54%  '

# Row 24: ImmutableX
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.14'
$ws.Range('E24').Value = '  -1.10%  '

# Row 26: PancakeSwap
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.43'
$ws.Range('E26').Value = '  -0.58%  '

# Row 27: LEO
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '4.03'
$ws.Range('E27').Value = '  -0.10%  '

# Row 28: EthereumClassic
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '25.57'
$ws.Range('E28').Value = '  +3.27%  '

# Row 29: Monero
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '165.86'
$ws.Range('E29').Value = '  -1.32%  '

# Row 30: Toncoin
$ws.Range('E30').Value = '  -0.21%  '

# Row 31: Cosmos
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '9.08'
$ws.Range('E31').Value = '  -0.60%  '

# Row 32: InjectiveProtocol
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '33.19'
$ws.Range('E32').Value = '  -0.33%  '

# Row 33: RenderToken
$ws.Range('E33').Value = '  +2.92%  '

# Row 34: FirstDigitalUSD
$ws.Range('E34').Value = '  +0.03%  '

# Row 35: Filecoin
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.05'
$ws.Range('E35').Value = '  -2.86%  '

# Row 36: Celestia
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '17.06'
$ws.Range('E36').Value = '  -6.23%  '

# Row 37: WEMIXToken
$ws.Range('E37').Value = '  -1.14%  '

# Row 38: Hedera
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0687'
$ws.Range('E38').Value = '  -0.70%  '

# Row 39: Kaspa
$ws.Range('E39').Value = '  -0.62%  '

# Row 40: ARBITRUM
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.77'
$ws.Range('E40').Value = '  -1.59%  '

# Row 41: LidoDAOToken
$ws.Range('B41').Value = 'LidoDAOToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.75'
$ws.Range('E41').Value = '  -0.35%  '

# Row 42: Stellar
$ws.Range('B42').Value = 'Stellar'
$ws.Range('C42').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.110'
$ws.Range('E42').Value = '  -1.54%  '

# Row 43: Maker
$ws.Range('D43').Value = '2.010.89'
$ws.Range('E43').Value = '  +0.68%  '

# Row 44: VeChain
$ws.Range('E44').Value = '  -2.33%  '

# Row 45: ApeXProtocol
$ws.Range('B45').Value = 'ApeXProtocol'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.17'
$ws.Range('E45').Value = '  +1.71%  '

# Row 46: FraxShare
$ws.Range('B46').Value = 'FraxShare'
$ws.Range('C46').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '10.14'
$ws.Range('E46').Value = '  +0.71%  '

# Row 47: EnergySwap
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '17.37'
$ws.Range('E47').Value = '  -0.39%  '

# Row 48: NEARProtocol
$ws.Range('E48').Value = '  -0.78%  '

# Row 49: MultiversX
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '53.85'
$ws.Range('E49').Value = '  -1.56%  '

# Row 50: RocketPoolETH
$ws.Range('D50').Value = '2.531.74'
$ws.Range('E50').Value = '  -0.12%  '

# Row 51: BitcoinSV
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '72.47'
$ws.Range('E51').Value = '  -1.16%  '
